# Adds 55 new leaderboard rows (game levels 1-3 were added, each level's
# completion appends a row of "<levelLetter>" + elapsed-time score to the
# leaderboard sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A: level identifier (stored as text); Column B: score/time value.
$levels = @("a", "a", "a", "a", "a", "l", "l", "a", "k", "a", "a", "l", "k", "k", "k", "k", "g", "a", "l", "a", "k", "a", "l", "k", "l", "a", "a", "a", "a", "l", "l", "l", "l", "a", "l", "l", "l", "a", "l", "l", "a", "a", "l", "l", "l", "a", "a", "a", "a", "a", "a", "l", "l", "a", "a")
$scores = @(1042.0, 1042.0, 1042.0, 1042.0, 521.0, 0.0, 1042.0, 1042.0, 1563.0, 1042.0, 1042.0, 1042.0, 1042.0, 0.0, 4689.0, 0.0, 2084.0, 1042.0, 1042.0, 1042.0, 1042.0, 0.0, 1042.0, 1042.0, 1563.0, 1042.0, 1042.0, 5210.0, 0.0, 0.0, 0.0, 3647.0, 2084.0, 1042.0, 0.0, 1042.0, 1042.0, 1563.0, 1042.0, 1042.0, 1042.0, 0.0, 0.0, 0.0, 0.0, 0.0, 2605.0, 0.0, 0.0, 1042.0, 3126.0, 1563.0, 2605.0, 0.0, 4168.0)

$startRow = 34
for ($i = 0; $i -lt $levels.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $levels[$i]
    $ws.Cells.Item($row, 2).Value = [double]$scores[$i]
}
